$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.825.25'
$ws.Range('E2').Value = '  +7.37%  '
$ws.Range('D3').Value = '2.420.62'
$ws.Range('E3').Value = '  +5.43%  '
$ws.Range('E4').Value = '  +0.69%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '116.85'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +12.57%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '319.65'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.59%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.637'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.04%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.633'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +4.67%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '43.37'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +10.22%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0940'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +4.71%  '
$ws.Range('E12').Value = '  +5.93%  '
$ws.Range('E13').Value = '  +4.70%  '
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.99'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.97%  '
$ws.Range('D16').Value = '2.791.36'
$ws.Range('D17').Value = '2.433.66'
$ws.Range('E17').Value = '  +6.10%  '
$ws.Range('D18').Value = '45.858.02'
$ws.Range('E18').Value = '  +7.59%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.64'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.56%  '
$ws.Range('E20').Value = '  +4.60%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.45'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '75.24'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.57'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.19'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.41'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +9.71%  '
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.67'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +6.55%  '
$ws.Range('E28').Value = '  +5.42%  '
$ws.Range('E29').Value = '  +0.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '40.36'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +11.69%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '23.07'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0973'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +13.94%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '173.73'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +5.48%  '
$ws.Range('E34').Value = '  +15.06%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.132'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.18%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.01'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +10.19%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.120'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +7.68%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.24'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +16.43%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.16'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +12.10%  '
$ws.Range('E40').Value = '  +5.94%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.82'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +15.45%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '102.06'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.80%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.71'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +13.09%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.240'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +5.82%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '72.75'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +2.81%  '
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.88'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +14.13%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '117.78'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.66%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.70'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +18.32%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '9.50'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +9.74%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '80.68'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.89%  '
